# Scheduled runner update: refresh Universalis-sourced market-price columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) on each job sheet's leve table.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 2448.0715
$ws.Range("I80").Value = 1727.2307
$ws.Range("J80").Value = 3072.8
$ws.Range("K80").Value = 5181.6921
$ws.Range("L80").Value = 9218.400000000001
$ws.Range("M80").Value = -4183.6921
$ws.Range("N80").Value = -11214.4

# Row 83
$ws.Range("H83").Value = 2448.0715
$ws.Range("I83").Value = 1727.2307
$ws.Range("J83").Value = 3072.8
$ws.Range("K83").Value = 15545.0763
$ws.Range("L83").Value = 27655.2
$ws.Range("M83").Value = -10553.0763
$ws.Range("N83").Value = -37639.2

# Row 132
$ws.Range("H132").Value = 2263.182
$ws.Range("I132").Value = 2144.8708
$ws.Range("K132").Value = 6434.6124
$ws.Range("M132").Value = -3904.6124

# Row 138
$ws.Range("H138").Value = 6539903.5
$ws.Range("I138").Value = 780.4737
$ws.Range("J138").Value = 10422508
$ws.Range("K138").Value = 2341.4211
$ws.Range("L138").Value = 31267524
$ws.Range("M138").Value = 2798.5789
$ws.Range("N138").Value = -31277804


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8596.290000000001
$ws.Range("I32").Value = 6025.5576
$ws.Range("K32").Value = 6025.5576
$ws.Range("M32").Value = -5738.5576


$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 14207.98
$ws.Range("I86").Value = 6724.7295
$ws.Range("J86").Value = 32666.666
$ws.Range("K86").Value = 6724.7295
$ws.Range("L86").Value = 32666.666
$ws.Range("M86").Value = -5601.7295
$ws.Range("N86").Value = -34912.666

# Row 89
$ws.Range("H89").Value = 14207.98
$ws.Range("I89").Value = 6724.7295
$ws.Range("J89").Value = 32666.666
$ws.Range("K89").Value = 33623.6475
$ws.Range("L89").Value = 163333.33
$ws.Range("M89").Value = -28007.6475
$ws.Range("N89").Value = -174565.33


$ws = $wb.Worksheets.Item("CRP")
# Row 5
$ws.Range("H5").Value = 2098.3333
$ws.Range("I5").Value = 3000
$ws.Range("J5").Value = 295
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 295
$ws.Range("M5").Value = -2888
$ws.Range("N5").Value = -519

# Row 14
$ws.Range("H14").Value = 800
$ws.Range("J14").Value = 800
$ws.Range("L14").Value = 800
$ws.Range("N14").Value = -1140

# Row 141
$ws.Range("H141").Value = 235998.2
$ws.Range("J141").Value = 235998.2
$ws.Range("L141").Value = 235998.2
$ws.Range("N141").Value = -246358.2


$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 7543.75
$ws.Range("I3").Value = 7087.5
$ws.Range("K3").Value = 21262.5
$ws.Range("M3").Value = -21150.5

# Row 37
$ws.Range("H37").Value = 157263.38
$ws.Range("J37").Value = 157263.38
$ws.Range("L37").Value = 471790.14
$ws.Range("N37").Value = -472014.14

# Row 68
$ws.Range("H68").Value = 2208.1052
$ws.Range("I68").Value = 2446
$ws.Range("J68").Value = 2180.1177
$ws.Range("K68").Value = 7338
$ws.Range("L68").Value = 6540.353099999999
$ws.Range("M68").Value = -6527
$ws.Range("N68").Value = -8162.353099999999

# Row 71
$ws.Range("H71").Value = 2208.1052
$ws.Range("I71").Value = 2446
$ws.Range("J71").Value = 2180.1177
$ws.Range("K71").Value = 22014
$ws.Range("L71").Value = 19621.0593
$ws.Range("M71").Value = -17958
$ws.Range("N71").Value = -27733.0593

# Row 134
$ws.Range("H134").Value = 1185.4286
$ws.Range("I134").Value = 1185.4286
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3556.2858
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 1513.7142
$ws.Range("N134").Value = ""

# Row 139
$ws.Range("H139").Value = 2032.2963
$ws.Range("I139").Value = 1866.9166
$ws.Range("K139").Value = 5600.7498
$ws.Range("M139").Value = -460.7497999999996

# Row 140
$ws.Range("H140").Value = 1172.1724
$ws.Range("I140").Value = 784.5
$ws.Range("J140").Value = 3033
$ws.Range("K140").Value = 2353.5
$ws.Range("L140").Value = 9099
$ws.Range("M140").Value = 2826.5
$ws.Range("N140").Value = -19459


$ws = $wb.Worksheets.Item("GSM")
# Row 12
$ws.Range("H12").Value = 4999.5
$ws.Range("I12").Value = 4999
$ws.Range("K12").Value = 4999
$ws.Range("M12").Value = -4859

# Row 21
$ws.Range("H21").Value = 19994
$ws.Range("I21").Value = 19994
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 19994
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -19821
$ws.Range("N21").Value = ""

# Row 29
$ws.Range("H29").Value = 6065.2
$ws.Range("I29").Value = 5079
$ws.Range("K29").Value = 5079
$ws.Range("M29").Value = -4789

# Row 30
$ws.Range("H30").Value = 19994
$ws.Range("I30").Value = 19994
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 19994
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -19889
$ws.Range("N30").Value = ""

# Row 80
$ws.Range("H80").Value = 2857.2
$ws.Range("I80").Value = 2649.6667
$ws.Range("K80").Value = 2649.6667
$ws.Range("M80").Value = -1651.6667

# Row 83
$ws.Range("H83").Value = 2857.2
$ws.Range("I83").Value = 2649.6667
$ws.Range("K83").Value = 13248.3335
$ws.Range("M83").Value = -8256.333500000001


$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2468.625
$ws.Range("I22").Value = 612.25
$ws.Range("J22").Value = 4325
$ws.Range("K22").Value = 612.25
$ws.Range("L22").Value = 4325
$ws.Range("M22").Value = -317.25
$ws.Range("N22").Value = -4915

# Row 27
$ws.Range("H27").Value = 2468.625
$ws.Range("I27").Value = 612.25
$ws.Range("J27").Value = 4325
$ws.Range("K27").Value = 612.25
$ws.Range("L27").Value = 4325
$ws.Range("M27").Value = -505.25
$ws.Range("N27").Value = -4539

# Row 68
$ws.Range("H68").Value = 2500
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 2000
$ws.Range("M68").Value = -1251

# Row 71
$ws.Range("H71").Value = 2500
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 10000
$ws.Range("M71").Value = -6256

# Row 136
$ws.Range("H136").Value = 2697.1333
$ws.Range("I136").Value = 689.625
$ws.Range("J136").Value = 4991.4287
$ws.Range("K136").Value = 2068.875
$ws.Range("L136").Value = 14974.2861
$ws.Range("M136").Value = 481.125
$ws.Range("N136").Value = -20074.2861

# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""


$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 8950
$ws.Range("I3").Value = 8950
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 8950
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -8836
$ws.Range("N3").Value = ""

# Row 31
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").Value = ""

# Row 81
$ws.Range("H81").Value = 1169.25
$ws.Range("I81").Value = 1033
$ws.Range("K81").Value = 2066
$ws.Range("M81").Value = -1005

# Row 84
$ws.Range("H84").Value = 1169.25
$ws.Range("I84").Value = 1033
$ws.Range("K84").Value = 10330
$ws.Range("M84").Value = -5026

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""

# Row 136
$ws.Range("H136").Value = 4323.1055
$ws.Range("I136").Value = 1271.2222
$ws.Range("J136").Value = 7069.8
$ws.Range("K136").Value = 3813.6666
$ws.Range("L136").Value = 21209.4
$ws.Range("M136").Value = -1263.6666
$ws.Range("N136").Value = -26309.4

# Row 137
$ws.Range("H137").Value = 121292.336
$ws.Range("J137").Value = 121292.336
$ws.Range("L137").Value = 121292.336
$ws.Range("N137").Value = -131492.336
